# Weekly update: insert two new price records (dated 44505) at the top of
# the "Vega Monumental Concepción - Pimiento" series, pushing the existing
# rows 168-200 down to 170-202.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 168:200 down by two rows, opening up
# a gap at rows 168:169 for the new entries.
$ws.Rows("168:169").Insert()

# New row 168 - Zafiro rojo
$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value = "Bíobío"
$ws.Cells.Item(168, 4).Value = 44505
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 100112002
$ws.Cells.Item(168, 7).Value = "Pimiento"
$ws.Cells.Item(168, 8).Value = "Zafiro rojo"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 180
$ws.Cells.Item(168, 11).Value = 42000
$ws.Cells.Item(168, 12).Value = 43000
$ws.Cells.Item(168, 13).Value = 42444
$ws.Cells.Item(168, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(168, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(168, 16).Value = 2830
$ws.Cells.Item(168, 17).Value = 15
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# New row 169 - Zafiro verde
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169, 3).Value = "Bíobío"
$ws.Cells.Item(169, 4).Value = 44505
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 100112002
$ws.Cells.Item(169, 7).Value = "Pimiento"
$ws.Cells.Item(169, 8).Value = "Zafiro verde"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 200
$ws.Cells.Item(169, 11).Value = 34000
$ws.Cells.Item(169, 12).Value = 35000
$ws.Cells.Item(169, 13).Value = 34500
$ws.Cells.Item(169, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(169, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(169, 16).Value = 2300
$ws.Cells.Item(169, 17).Value = 15
$ws.Cells.Item(169, 18).Value = "Hortaliza"
